{"js": "// Remove the empty paragraph, the page-break paragraph, and the\n// \"\u00a9 2020 ...\" copyright paragraph that used to follow the\n// \"LOB1011: Eletricidade Aplicada (Requisito fraco)\" paragraph.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  items[i].load(\"text\");\n}\nawait context.sync();\n\n// Locate the anchor paragraph by its text.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOB1011\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find the 'LOB1011' anchor paragraph.\");\n}\n\n// Delete the three paragraphs that immediately follow the anchor\n// (blank paragraph, page-break paragraph, and the copyright paragraph).\nconst toDelete = [];\nfor (let i = anchorIndex + 1; i <= anchorIndex + 3 && i < items.length; i++) {\n  toDelete.push(items[i]);\n}\n\nfor (const para of toDelete) {\n  para.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the empty paragraph, the page-break paragraph, and the\n# \"\u00a9 2020 ...\" copyright paragraph that used to follow the\n# \"LOB1011: Eletricidade Aplicada (Requisito fraco)\" paragraph.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph by its text.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*LOB1011*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find the 'LOB1011' anchor paragraph.\"\n}\n\n# Delete the three paragraphs that immediately follow the anchor\n# (blank paragraph, page-break paragraph, and the copyright paragraph).\n# Deleting repeatedly at (anchorIndex + 1) works because each deletion\n# shifts the following paragraphs up by one.\nfor ($n = 0; $n -lt 3; $n++) {\n    $p = $d.Paragraphs.Item($anchorIndex + 1)\n    $p.Range.Delete()\n}\n"}
